$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '50.952.71'
$ws.Cells.Item(2, 5).Value = '  -0.50%  '
$ws.Cells.Item(3, 4).Value = '2.948.45'
$ws.Cells.Item(3, 5).Value = '  -0.25%  '
$ws.Cells.Item(4, 5).Value = '  -0.02%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '378.76'
$ws.Cells.Item(5, 5).Value = '  -1.00%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '101.52'
$ws.Cells.Item(6, 5).Value = '  -1.16%  '
$ws.Cells.Item(7, 5).Value = '  +0.44%  '
$ws.Cells.Item(8, 5).Value = '  -0.03%  '
$ws.Cells.Item(9, 5).Value = '  -1.36%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '36.18'
$ws.Cells.Item(10, 5).Value = '  -1.30%  '
$ws.Cells.Item(11, 5).Value = '  -0.52%  '
$ws.Cells.Item(12, 5).Value = '  +0.75%  '
$ws.Cells.Item(13, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(13, 4).Value = '3.413.15'
$ws.Cells.Item(13, 5).Value = '  -0.42%  '
$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '7.78'
$ws.Cells.Item(14, 5).Value = '  +4.75%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '18.29'
$ws.Cells.Item(15, 5).Value = '  +1.11%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '12.09'
$ws.Cells.Item(16, 5).Value = '  +68.89%  '
$ws.Cells.Item(17, 4).Value = '2.951.17'
$ws.Cells.Item(17, 5).Value = '  -0.36%  '
$ws.Cells.Item(18, 5).Value = '  +1.81%  '
$ws.Cells.Item(19, 4).Value = '50.915.02'
$ws.Cells.Item(19, 5).Value = '  -0.56%  '
$ws.Cells.Item(21, 5).Value = '  -1.67%  '
$ws.Cells.Item(22, 4).Value = '0.0₃0952'
$ws.Cells.Item(22, 5).Value = '  -0.43%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '69.33'
$ws.Cells.Item(23, 5).Value = '  +1.07%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '266.13'
$ws.Cells.Item(24, 5).Value = '  +1.36%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '3.18'
$ws.Cells.Item(25, 5).Value = '  +8.64%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '8.10'
$ws.Cells.Item(26, 5).Value = '  -3.08%  '
$ws.Cells.Item(27, 5).Value = '  +0.00%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '7.03'
$ws.Cells.Item(28, 5).Value = '  -8.10%  '
$ws.Cells.Item(30, 5).Value = '  -3.72%  '
$ws.Cells.Item(31, 5).Value = '  -3.84%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '10.18'
$ws.Cells.Item(32, 5).Value = '  +3.47%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '50.52'
$ws.Cells.Item(33, 5).Value = '  -0.10%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '2.06'
$ws.Cells.Item(34, 5).Value = '  +0.01%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '33.57'
$ws.Cells.Item(35, 5).Value = '  -2.15%  '
$ws.Cells.Item(36, 5).Value = '  -5.30%  '
$ws.Cells.Item(37, 5).Value = '  -0.06%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '3.10'
$ws.Cells.Item(38, 5).Value = '  +3.76%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '16.57'
$ws.Cells.Item(40, 5).Value = '  -1.88%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.82'
$ws.Cells.Item(41, 5).Value = '  +1.42%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '2.51'
$ws.Cells.Item(42, 5).Value = '  -2.35%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '119.43'
$ws.Cells.Item(43, 5).Value = '  -1.42%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '3.52'
$ws.Cells.Item(44, 5).Value = '  +8.29%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '21.39'
$ws.Cells.Item(45, 5).Value = '  -0.36%  '
$ws.Cells.Item(46, 5).Value = '  -1.96%  '
$ws.Cells.Item(47, 5).Value = '  -2.42%  '
$ws.Cells.Item(48, 4).Value = '2.000.70'
$ws.Cells.Item(48, 5).Value = '  -0.72%  '
$ws.Cells.Item(49, 5).Value = '  -4.48%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.0317'
$ws.Cells.Item(50, 5).Value = '  -9.73%  '
$ws.Cells.Item(51, 5).Value = '  +3.80%  '
